# Increased Slugs and Buckshot damages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slugs (row 38): raw damage value H38 2.5 -> 2.7
$ws.Range("H38").Value = 2.7

# Buckshot (row 39): formula H39 "=9*0.4" -> "=9*0.42"
$ws.Range("H39").Formula = "=9*0.42"

# Move the selection cursor to match the saved view state (N22 -> J27)
$ws.Range("J27").Select()
